# Refresh the crypto price/volume snapshot cells to match the latest scrape.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'27.696.62"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = "`'1.895.45"
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  -1.13%  '
$ws.Range("D5").Value = "`'311.86"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("D7").Value = "`'0.4888"
$ws.Range("E7").Value = '  +1.21%  '
$ws.Range("D8").Value = "`'0.3789"
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").Value = "`'0.9127"
$ws.Range("E10").Value = '  -2.88%  '
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "`'1.959.41"
$ws.Range("E12").Value = '  +4.38%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "`'0.07663"
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").Value = "`'5.481"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = "`'6.602"
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = "`'91.25"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = "`'0.000008763"
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = "`'27.570.87"
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("E21").Value = '  -2.44%  '
$ws.Range("D22").Value = "`'5.120"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = "`'2.095.58"
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("D24").Value = "`'10.75"
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").Value = "`'1.905"
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").Value = "`'153.73"
$ws.Range("E26").Value = '  -2.52%  '
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").Value = "`'2.152"
$ws.Range("E28").Value = '  +5.31%  '
$ws.Range("D29").Value = "`'115.46"
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = "`'4.873"
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").Value = '  -4.25%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = "`'0.7658"
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("D35").Value = "`'4.632"
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = "`'0.02036"
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = "`'2.534"
$ws.Range("E37").Value = '  -7.05%  '
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = "`'0.5472"
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("D41").Value = "`'2.977"
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("D42").Value = "`'6.879"
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("D43").Value = "`'8.538"
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = "`'0.1517"
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").Value = "`'112.13"
$ws.Range("E45").Value = '  +6.35%  '
$ws.Range("D46").Value = "`'10.68"
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").Value = "`'0.4782"
$ws.Range("E47").Value = '  -1.84%  '
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").Value = "`'1.640"
$ws.Range("D50").Value = "`'67.45"
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").Value = "`'0.06050"
$ws.Range("E51").Value = '  -1.23%  '
